# "added few more programs for evaluation and minor bug fixes"
# Rows 1-20 are unchanged. Starting at row 21, two new p02879 entries are
# inserted (pushing the former rows 21-33 down by two), then a further
# block of new program rows is appended, ending with a relocated p03795
# row and a brand-new trailing p03544 row. Net effect: the sheet grows
# from A1:E33 to A1:E47, so every row from 21 onward is simply rewritten
# in place with its final target content (values only; no formulas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(21, 1).Value = "p02879"
$ws.Cells.Item(21, 2).Value = "s055491155.cob"
$ws.Cells.Item(21, 3).Value = "s055491155.py"
$ws.Cells.Item(21, 4).Value = 100
$ws.Cells.Item(21, 5).Value = "23/23"

$ws.Cells.Item(22, 1).Value = "p02879"
$ws.Cells.Item(22, 2).Value = "s914345468.cob"
$ws.Cells.Item(22, 3).Value = "s914345468.py"
$ws.Cells.Item(22, 4).Value = 100
$ws.Cells.Item(22, 5).Value = "23/23"

$ws.Cells.Item(23, 1).Value = "p02909"
$ws.Cells.Item(23, 2).Value = "s502766942.cob"
$ws.Cells.Item(23, 3).Value = "s502766942.py"
$ws.Cells.Item(23, 4).Value = 100
$ws.Cells.Item(23, 5).Value = "3/3"

$ws.Cells.Item(24, 1).Value = "p02909"
$ws.Cells.Item(24, 2).Value = "s880642132.cob"
$ws.Cells.Item(24, 3).Value = "s880642132.py"
$ws.Cells.Item(24, 4).Value = 100
$ws.Cells.Item(24, 5).Value = "3/3"

$ws.Cells.Item(25, 1).Value = "p02915"
$ws.Cells.Item(25, 2).Value = "s250806848.cob"
$ws.Cells.Item(25, 3).Value = "s250806848.py"
$ws.Cells.Item(25, 4).Value = 100
$ws.Cells.Item(25, 5).Value = "9/9"

$ws.Cells.Item(26, 1).Value = "p02993"
$ws.Cells.Item(26, 2).Value = "s111459231.cob"
$ws.Cells.Item(26, 3).Value = "s111459231.py"
$ws.Cells.Item(26, 4).Value = 100
$ws.Cells.Item(26, 5).Value = "13/13"

$ws.Cells.Item(27, 1).Value = "p03029"
$ws.Cells.Item(27, 2).Value = "s018605057.cob"
$ws.Cells.Item(27, 3).Value = "s018605057.py"
$ws.Cells.Item(27, 4).Value = 100
$ws.Cells.Item(27, 5).Value = "11/11"

$ws.Cells.Item(28, 1).Value = "p03085"
$ws.Cells.Item(28, 2).Value = "s247348869.cob"
$ws.Cells.Item(28, 3).Value = "s247348869.py"
$ws.Cells.Item(28, 4).Value = 100
$ws.Cells.Item(28, 5).Value = "4/4"

$ws.Cells.Item(29, 1).Value = "p03085"
$ws.Cells.Item(29, 2).Value = "s794500343.cob"
$ws.Cells.Item(29, 3).Value = "s794500343.py"
$ws.Cells.Item(29, 4).Value = 100
$ws.Cells.Item(29, 5).Value = "4/4"

$ws.Cells.Item(30, 1).Value = "p03101"
$ws.Cells.Item(30, 2).Value = "s563515874.cob"
$ws.Cells.Item(30, 3).Value = "s563515874.py"
$ws.Cells.Item(30, 4).Value = 100
$ws.Cells.Item(30, 5).Value = "15/15"

$ws.Cells.Item(31, 1).Value = "p03238"
$ws.Cells.Item(31, 2).Value = "s739597451.cob"
$ws.Cells.Item(31, 3).Value = "s739597451.py"
$ws.Cells.Item(31, 4).Value = 100
$ws.Cells.Item(31, 5).Value = "7/7"

$ws.Cells.Item(32, 1).Value = "p03260"
$ws.Cells.Item(32, 2).Value = "s346414249.cob"
$ws.Cells.Item(32, 3).Value = "s346414249.py"
$ws.Cells.Item(32, 4).Value = 100
$ws.Cells.Item(32, 5).Value = "9/9"

$ws.Cells.Item(33, 1).Value = "p03315"
$ws.Cells.Item(33, 2).Value = "s910259082.cob"
$ws.Cells.Item(33, 3).Value = "s910259082.py"
$ws.Cells.Item(33, 4).Value = 100
$ws.Cells.Item(33, 5).Value = "8/8"

$ws.Cells.Item(34, 1).Value = "p03316"
$ws.Cells.Item(34, 2).Value = "s513397080.cob"
$ws.Cells.Item(34, 3).Value = "s513397080.py"
$ws.Cells.Item(34, 4).Value = 100
$ws.Cells.Item(34, 5).Value = "11/11"

$ws.Cells.Item(35, 1).Value = "p03331"
$ws.Cells.Item(35, 2).Value = "s251999208.cob"
$ws.Cells.Item(35, 3).Value = "s251999208.py"
$ws.Cells.Item(35, 4).Value = 50
$ws.Cells.Item(35, 5).Value = "1/2"

$ws.Cells.Item(36, 1).Value = "p03415"
$ws.Cells.Item(36, 2).Value = "s342220458.cob"
$ws.Cells.Item(36, 3).Value = "s342220458.py"
$ws.Cells.Item(36, 4).Value = 100
$ws.Cells.Item(36, 5).Value = "2/2"

$ws.Cells.Item(37, 1).Value = "p03415"
$ws.Cells.Item(37, 2).Value = "s760213038.cob"
$ws.Cells.Item(37, 3).Value = "s760213038.py"
$ws.Cells.Item(37, 4).Value = 100
$ws.Cells.Item(37, 5).Value = "2/2"

$ws.Cells.Item(38, 1).Value = "p03433"
$ws.Cells.Item(38, 2).Value = "s378680164.cob"
$ws.Cells.Item(38, 3).Value = "s378680164.py"
$ws.Cells.Item(38, 4).Value = 100
$ws.Cells.Item(38, 5).Value = "3/3"

$ws.Cells.Item(39, 1).Value = "p03433"
$ws.Cells.Item(39, 2).Value = "s910105267.cob"
$ws.Cells.Item(39, 3).Value = "s910105267.py"
$ws.Cells.Item(39, 4).Value = 100
$ws.Cells.Item(39, 5).Value = "3/3"

$ws.Cells.Item(40, 1).Value = "p03493"
$ws.Cells.Item(40, 2).Value = "s535106378.cob"
$ws.Cells.Item(40, 3).Value = "s535106378.py"
$ws.Cells.Item(40, 4).Value = 100
$ws.Cells.Item(40, 5).Value = "2/2"

$ws.Cells.Item(41, 1).Value = "p03605"
$ws.Cells.Item(41, 2).Value = "s655098455.cob"
$ws.Cells.Item(41, 3).Value = "s655098455.py"
$ws.Cells.Item(41, 4).Value = 100
$ws.Cells.Item(41, 5).Value = "90/90"

$ws.Cells.Item(42, 1).Value = "p03605"
$ws.Cells.Item(42, 2).Value = "s724623140.cob"
$ws.Cells.Item(42, 3).Value = "s724623140.py"
$ws.Cells.Item(42, 4).Value = 100
$ws.Cells.Item(42, 5).Value = "90/90"

$ws.Cells.Item(43, 1).Value = "p03623"
$ws.Cells.Item(43, 2).Value = "s498016040.cob"
$ws.Cells.Item(43, 3).Value = "s498016040.py"
$ws.Cells.Item(43, 4).Value = 100
$ws.Cells.Item(43, 5).Value = "2/2"

$ws.Cells.Item(44, 1).Value = "p03693"
$ws.Cells.Item(44, 2).Value = "s862197544.cob"
$ws.Cells.Item(44, 3).Value = "s862197544.py"
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = "0/2"

$ws.Cells.Item(45, 1).Value = "p03737"
$ws.Cells.Item(45, 2).Value = "s496684777.cob"
$ws.Cells.Item(45, 3).Value = "s496684777.py"
$ws.Cells.Item(45, 4).Value = 100
$ws.Cells.Item(45, 5).Value = "4/4"

$ws.Cells.Item(46, 1).Value = "p03795"
$ws.Cells.Item(46, 2).Value = "s725157986.cob"
$ws.Cells.Item(46, 3).Value = "s725157986.py"
$ws.Cells.Item(46, 4).Value = 100
$ws.Cells.Item(46, 5).Value = "2/2"

$ws.Cells.Item(47, 1).Value = "p03544"
$ws.Cells.Item(47, 2).Value = "s498531048.cob"
$ws.Cells.Item(47, 3).Value = "s498531048.py"
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = "0/13"

